$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'vgr_2'
$ws.Cells.Item(2, 2).Value = 885
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = '/vgr/pick_up_and_transport'
$ws.Cells.Item(2, 5).Value = '4d198444-6633-4218-b1f7-ca67ec666360'
$ws.Cells.Item(2, 6).Value = '[''parameter_start_position'', ''parameter_end_position'']'

$ws.Cells.Item(3, 1).Value = 'mm_2'
$ws.Cells.Item(3, 2).Value = 381
$ws.Cells.Item(3, 3).Value = 4
$ws.Cells.Item(3, 4).Value = '/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to'
$ws.Cells.Item(3, 5).Value = '570d0814-988a-4856-bc82-249db6050f5e'
$ws.Cells.Item(3, 6).Value = '[''parameter_start_position'', ''parameter_end_position'', ''parameter_burn_workpiece_size'', ''parameter_quantity'']'

$ws.Cells.Item(4, 1).Value = 'sm_1'
$ws.Cells.Item(4, 2).Value = 378
$ws.Cells.Item(4, 3).Value = 2
$ws.Cells.Item(4, 4).Value = '/sm/sort, /sm/transport'
$ws.Cells.Item(4, 5).Value = '16d2bd16-3be9-4daa-a4ad-edb7f5818fcb'
$ws.Cells.Item(4, 6).Value = '[''parameter_use_nfc'', ''parameter_start_position'', ''parameter_end_position'', ''parameter_sorting_machine_ejection_position'']'

$ws.Cells.Item(5, 1).Value = 'hbw_1'
$ws.Cells.Item(5, 2).Value = 873
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = '/hbw/get_empty_bucket, /hbw/store, /hbw/store_empty_bucket, /hbw/unload'
$ws.Cells.Item(5, 5).Value = 'b179f074-238d-4666-b50f-9a8959d0a48e'
$ws.Cells.Item(5, 6).Value = '[''parameter_hbw_slot'', ''parameter_use_nfc'']'

$ws.Cells.Item(6, 1).Value = 'wt_1'
$ws.Cells.Item(6, 2).Value = 447
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = '/wt/pick_up_and_transport'
$ws.Cells.Item(6, 5).Value = '8febb390-19ce-4d63-a018-d9617a8bb1b7'
$ws.Cells.Item(6, 6).Value = '[''parameter_start_position'', ''parameter_end_position'']'

$ws.Cells.Item(7, 1).Value = 'sm_2'
$ws.Cells.Item(7, 2).Value = 309
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = '/sm/sort, /sm/transport'
$ws.Cells.Item(7, 5).Value = '722f5091-ed89-45a3-89c7-4962901b6c14'
$ws.Cells.Item(7, 6).Value = '[''parameter_start_position'', ''parameter_end_position'']'

$ws.Cells.Item(8, 1).Value = 'vgr_1'
$ws.Cells.Item(8, 2).Value = 1866
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = '/vgr/pick_up_and_transport'
$ws.Cells.Item(8, 5).Value = '0e7b5a4c-4c03-47b2-96fd-e401ed7fbca9'
$ws.Cells.Item(8, 6).Value = '[''parameter_start_position'', ''parameter_end_position'']'

$ws.Cells.Item(9, 1).Value = 'ov_1'
$ws.Cells.Item(9, 2).Value = 612
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 4).Value = '/ov/burn, /ov/temper'
$ws.Cells.Item(9, 5).Value = '633d065f-96c0-4c4b-8112-302990575763'
$ws.Cells.Item(9, 6).Value = '[''parameter_burn_workpiece_size'', ''parameter_burn_workpiece_thickness'']'

$ws.Cells.Item(10, 1).Value = 'pm_1'
$ws.Cells.Item(10, 2).Value = 204
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = '/pm/punch_gill, /pm/punch_recesses, /pm/punch_ribbing'
$ws.Cells.Item(10, 5).Value = '21559c95-22a5-4c8b-9424-dbbc14a9f63b'
$ws.Cells.Item(10, 6).Value = '[''parameter_start_position'', ''parameter_end_position'', ''parameter_quantity'']'

$ws.Cells.Item(11, 1).Value = 'mm_1'
$ws.Cells.Item(11, 2).Value = 576
$ws.Cells.Item(11, 3).Value = 4
$ws.Cells.Item(11, 4).Value = '/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to'
$ws.Cells.Item(11, 5).Value = '167db95e-ae8b-4ae8-ac11-055401e11894'
$ws.Cells.Item(11, 6).Value = '[''parameter_start_position'', ''parameter_end_position'']'

$ws.Cells.Item(12, 1).Value = 'hbw_2'
$ws.Cells.Item(12, 2).Value = 1581
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = '/hbw/store_empty_bucket, /hbw/unload'
$ws.Cells.Item(12, 5).Value = '27958fc0-4484-41ff-9260-e76f8a83a7cd'
$ws.Cells.Item(12, 6).Value = '[''parameter_hbw_slot'', ''parameter_use_nfc'']'

$ws.Cells.Item(13, 1).Value = 'ov_2'
$ws.Cells.Item(13, 2).Value = 330
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = '/ov/burn'
$ws.Cells.Item(13, 5).Value = '1ab1350e-cba4-42ea-8efd-a0b01e88380e'
$ws.Cells.Item(13, 6).Value = '[''parameter_burn_workpiece_size'', ''parameter_burn_workpiece_thickness'']'

$ws.Cells.Item(14, 1).Value = 'wt_2'
$ws.Cells.Item(14, 2).Value = 330
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = '/wt/pick_up_and_transport'
$ws.Cells.Item(14, 5).Value = '7316381c-127f-43cb-956b-ca72e60bc6ab'
$ws.Cells.Item(14, 6).Value = '[''parameter_start_position'', ''parameter_end_position'']'

$ws.Cells.Item(15, 1).Value = 'dm_2'
$ws.Cells.Item(15, 2).Value = 177
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 4).Value = '/dm/cylindrical_drill, /dm/drill, /dm/lower'
$ws.Cells.Item(15, 5).Value = 'ad6c9c0b-f3ba-45e7-b887-b96bf0260887'
$ws.Cells.Item(15, 6).Value = '[''parameter_start_position'', ''parameter_end_position'']'

$ws.Cells.Item(16, 1).Value = 'hw_1'
$ws.Cells.Item(16, 2).Value = 522
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = '/hw/human_review'
$ws.Cells.Item(16, 5).Value = 'a8d0fcdd-46c6-44f8-8b19-cecd803d356f'
$ws.Cells.Item(16, 6).Value = '[]'
